# Auto-generated Excel COM-interop script applying market-data refresh
# to the Chocobo_Profits workbook (H:N columns per Leve row) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 4
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 105.28571
$ws.Range("I4").Value = 106.166664
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 106.166664
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 7.833336000000003
$ws.Range("N4").Value = -328

# Sheet ALC, row 68
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 32500
$ws.Range("J68").Value = 32500
$ws.Range("L68").Value = 32500
$ws.Range("N68").Value = -33998

# Sheet ALC, row 71
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 32500
$ws.Range("J71").Value = 32500
$ws.Range("L71").Value = 97500
$ws.Range("N71").Value = -104988

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4020.3
$ws.Range("I86").Value = 2267.6667
$ws.Range("J86").Value = 4771.4287
$ws.Range("K86").Value = 2267.6667
$ws.Range("L86").Value = 4771.4287
$ws.Range("M86").Value = -1144.6667
$ws.Range("N86").Value = -7017.4287

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4020.3
$ws.Range("I89").Value = 2267.6667
$ws.Range("J89").Value = 4771.4287
$ws.Range("K89").Value = 11338.3335
$ws.Range("L89").Value = 23857.1435
$ws.Range("M89").Value = -5722.333500000001
$ws.Range("N89").Value = -35089.14350000001

# Sheet ALC, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1032.381
$ws.Range("J129").Value = 1093.8948
$ws.Range("L129").Value = 3281.6844
$ws.Range("N129").Value = -13281.6844

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 28576886
$ws.Range("I132").Value = 33338796
$ws.Range("J132").Value = 5416
$ws.Range("K132").Value = 100016388
$ws.Range("L132").Value = 16248
$ws.Range("M132").Value = -100013858
$ws.Range("N132").Value = -21308

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3260.4285
$ws.Range("I137").Value = 3189.6
$ws.Range("J137").Value = 3378.476
$ws.Range("K137").Value = 9568.799999999999
$ws.Range("L137").Value = 10135.428
$ws.Range("M137").Value = -7018.799999999999
$ws.Range("N137").Value = -15235.428

# Sheet ARM, row 62
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Sheet ARM, row 65
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3311.647
$ws.Range("I122").Value = 1729.9
$ws.Range("K122").Value = 5189.700000000001
$ws.Range("M122").Value = -2739.700000000001

# Sheet BSM, row 16
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 3000
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3340

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4758.0527
$ws.Range("I31").Value = 2120.3
$ws.Range("J31").Value = 7688.8887
$ws.Range("K31").Value = 2120.3
$ws.Range("L31").Value = 7688.8887
$ws.Range("M31").Value = -1825.3
$ws.Range("N31").Value = -8278.8887

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4758.0527
$ws.Range("I34").Value = 2120.3
$ws.Range("J34").Value = 7688.8887
$ws.Range("K34").Value = 2120.3
$ws.Range("L34").Value = 7688.8887
$ws.Range("M34").Value = -1918.3
$ws.Range("N34").Value = -8092.8887

# Sheet CRP, row 63
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 42496.9
$ws.Range("I63").Value = 9999
$ws.Range("J63").Value = 46107.777
$ws.Range("K63").Value = 9999
$ws.Range("L63").Value = 46107.777
$ws.Range("M63").Value = -9313
$ws.Range("N63").Value = -47479.777

# Sheet CRP, row 66
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 42496.9
$ws.Range("I66").Value = 9999
$ws.Range("J66").Value = 46107.777
$ws.Range("K66").Value = 29997
$ws.Range("L66").Value = 138323.331
$ws.Range("M66").Value = -26565
$ws.Range("N66").Value = -145187.331

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2311.0833
$ws.Range("I5").Value = 693.6
$ws.Range("K5").Value = 2080.8
$ws.Range("M5").Value = -1968.8

# Sheet CUL, row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3846.6
$ws.Range("I80").Value = 1850
$ws.Range("J80").Value = 4572.636
$ws.Range("K80").Value = 5550
$ws.Range("L80").Value = 13717.908
$ws.Range("M80").Value = -4614
$ws.Range("N80").Value = -15589.908

# Sheet CUL, row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3846.6
$ws.Range("I83").Value = 1850
$ws.Range("J83").Value = 4572.636
$ws.Range("K83").Value = 16650
$ws.Range("L83").Value = 41153.724
$ws.Range("M83").Value = -11970
$ws.Range("N83").Value = -50513.724

# Sheet CUL, row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2083.1865
$ws.Range("I121").Value = 437.5
$ws.Range("J121").Value = 2202.8728
$ws.Range("K121").Value = 1312.5
$ws.Range("L121").Value = 6608.6184
$ws.Range("M121").Value = -2.5
$ws.Range("N121").Value = -9228.618399999999

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2940.7886
$ws.Range("I122").Value = 723.3
$ws.Range("J122").Value = 3468.762
$ws.Range("K122").Value = 6509.7
$ws.Range("L122").Value = 31218.858
$ws.Range("M122").Value = -4059.7
$ws.Range("N122").Value = -36118.858

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7940668
$ws.Range("J131").Value = 1015.8409
$ws.Range("L131").Value = 3047.5227
$ws.Range("N131").Value = -13127.5227

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2311.0833
$ws.Range("I135").Value = 693.6
$ws.Range("K135").Value = 6242.400000000001
$ws.Range("M135").Value = -3707.400000000001

# Sheet GSM, row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29831.666
$ws.Range("I57").Value = 32666.666
$ws.Range("J57").Value = 28886.666
$ws.Range("K57").Value = 32666.666
$ws.Range("L57").Value = 28886.666
$ws.Range("M57").Value = -31846.666
$ws.Range("N57").Value = -30526.666

# Sheet GSM, row 68
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 40295
$ws.Range("J68").Value = 40295
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41917

# Sheet GSM, row 71
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H71").Value = 40295
$ws.Range("J71").Value = 40295
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128997

# Sheet LTW, row 62
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 39624.5
$ws.Range("J62").Value = 39624.5
$ws.Range("L62").Value = 39624.5
$ws.Range("N62").Value = -40872.5

# Sheet LTW, row 65
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 39624.5
$ws.Range("J65").Value = 39624.5
$ws.Range("L65").Value = 118873.5
$ws.Range("N65").Value = -125113.5

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5347.96
$ws.Range("I132").Value = 1205.3158
$ws.Range("K132").Value = 3615.9474
$ws.Range("M132").Value = -1085.9474

# Sheet LTW, row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 51530.11
$ws.Range("J139").Value = 51530.11
$ws.Range("L139").Value = 51530.11
$ws.Range("N139").Value = -61810.11

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3614.5386
$ws.Range("I122").Value = 2121
$ws.Range("K122").Value = 6363
$ws.Range("M122").Value = -3913

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9013499
$ws.Range("I132").Value = 4852.4585
$ws.Range("J132").Value = 25644848
$ws.Range("K132").Value = 14557.3755
$ws.Range("L132").Value = 76934544
$ws.Range("M132").Value = -12027.3755
$ws.Range("N132").Value = -76939604

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3873.3845
$ws.Range("I136").Value = 1076.1538
$ws.Range("J136").Value = 6670.615
$ws.Range("K136").Value = 3228.4614
$ws.Range("L136").Value = 20011.845
$ws.Range("M136").Value = -678.4614000000001
$ws.Range("N136").Value = -25111.845
